$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2634.6316
$ws.Range("J17").Value = 2634.6316
$ws.Range("L17").Value = 7903.8948
$ws.Range("N17").Value = -8239.8948

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()

$ws.Range("H116").Value = 3749.25
$ws.Range("I116").Value = 2499
$ws.Range("K116").Value = 2499
$ws.Range("M116").Value = 943

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3799.25
$ws.Range("I2").Value = 3506.9092
$ws.Range("K2").Value = 3506.9092
$ws.Range("M2").Value = -3393.9092

$ws.Range("H4").Value = 999
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1232

$ws.Range("H97").Value = 1130.5
$ws.Range("I97").Value = 887.9
$ws.Range("K97").Value = 887.9
$ws.Range("M97").Value = -391.9

$ws.Range("H116").Value = 3799.25
$ws.Range("I116").Value = 3506.9092
$ws.Range("K116").Value = 3506.9092
$ws.Range("M116").Value = -1212.9092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3799.25
$ws.Range("I3").Value = 3506.9092
$ws.Range("K3").Value = 3506.9092
$ws.Range("M3").Value = -3392.9092

$ws.Range("H140").Value = 141712
$ws.Range("J140").Value = 141712
$ws.Range("L140").Value = 141712
$ws.Range("N140").Value = -152072

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1471.7
$ws.Range("I16").Value = 1403.8
$ws.Range("J16").Value = 1539.6
$ws.Range("K16").Value = 1403.8
$ws.Range("L16").Value = 1539.6
$ws.Range("M16").Value = -1116.8
$ws.Range("N16").Value = -2113.6

$ws.Range("H31").Value = 5998.0835
$ws.Range("I31").Value = 3206.4167
$ws.Range("K31").Value = 3206.4167
$ws.Range("M31").Value = -2911.4167

$ws.Range("H34").Value = 5998.0835
$ws.Range("I34").Value = 3206.4167
$ws.Range("K34").Value = 3206.4167
$ws.Range("M34").Value = -3004.4167

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H58").Value = 4356.9
$ws.Range("I58").Value = 2248.1667
$ws.Range("K58").Value = 2248.1667
$ws.Range("M58").Value = -2045.1667

$ws.Range("H74").Value = 60500.355
$ws.Range("J74").Value = 60500.355
$ws.Range("L74").Value = 60500.355
$ws.Range("N74").Value = -62248.355

$ws.Range("H77").Value = 60500.355
$ws.Range("J77").Value = 60500.355
$ws.Range("L77").Value = 181501.065
$ws.Range("N77").Value = -190237.065

$ws.Range("H105").Value = 1681.25
$ws.Range("I105").Value = 2412.5
$ws.Range("J105").Value = 950
$ws.Range("K105").Value = 2412.5
$ws.Range("L105").Value = 950
$ws.Range("M105").Value = -665.5
$ws.Range("N105").Value = -4444

$ws.Range("H113").Value = 1471.7
$ws.Range("I113").Value = 1403.8
$ws.Range("J113").Value = 1539.6
$ws.Range("K113").Value = 1403.8
$ws.Range("L113").Value = 1539.6
$ws.Range("M113").Value = 766.2
$ws.Range("N113").Value = -5879.6

$ws.Range("H122").Value = 1272.4286
$ws.Range("I122").Value = 1201.6
$ws.Range("K122").Value = 3604.8
$ws.Range("M122").Value = -1154.8

$ws.Range("H136").Value = 4356.9
$ws.Range("I136").Value = 2248.1667
$ws.Range("K136").Value = 6744.500100000001
$ws.Range("M136").Value = -4194.500100000001

$ws.Range("H141").Value = 80624.10000000001
$ws.Range("J141").Value = 80624.10000000001
$ws.Range("L141").Value = 80624.10000000001
$ws.Range("N141").Value = -90984.10000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 142914.72
$ws.Range("I4").Value = 333348.66
$ws.Range("K4").Value = 1000045.98
$ws.Range("M4").Value = -999933.98

$ws.Range("H12").Value = 143.83333
$ws.Range("J12").Value = 192.25
$ws.Range("L12").Value = 576.75
$ws.Range("N12").Value = -922.75

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1792.625
$ws.Range("I102").Value = 1792.625
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1792.625
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -170.625
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 13000
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 25000
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 25000
$ws.Range("M3").Value = -888
$ws.Range("N3").Value = -25224

$ws.Range("H9").Value = 289.6
$ws.Range("I9").Value = 150
$ws.Range("J9").Value = 499
$ws.Range("K9").Value = 150
$ws.Range("L9").Value = 499
$ws.Range("M9").Value = 74
$ws.Range("N9").Value = -947

$ws.Range("H10").Value = 1097.6
$ws.Range("J10").Value = 122.25
$ws.Range("L10").Value = 122.25
$ws.Range("N10").Value = -402.25

$ws.Range("H11").Value = 833.3333
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = -360
$ws.Range("N11").Value = -1280

$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("N13").Value = -1280

$ws.Range("H15").Value = 13000
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 25000
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 25000
$ws.Range("M15").Value = -830
$ws.Range("N15").Value = -25340

$ws.Range("H17").Value = 4000
$ws.Range("I17").Value = 2000
$ws.Range("K17").Value = 2000
$ws.Range("M17").Value = -1830

$ws.Range("H19").Value = 5000
$ws.Range("I19").Value = 5000
$ws.Range("K19").Value = 5000
$ws.Range("M19").Value = -4830

$ws.Range("H25").Value = 3625
$ws.Range("J25").Value = 4000
$ws.Range("L25").Value = 4000
$ws.Range("N25").Value = -4460

$ws.Range("H46").Value = 5471.5293
$ws.Range("I46").Value = 4503.5
$ws.Range("J46").Value = 5999.5454
$ws.Range("K46").Value = 4503.5
$ws.Range("L46").Value = 5999.5454
$ws.Range("M46").Value = -4315.5
$ws.Range("N46").Value = -6375.5454

$ws.Range("H122").Value = 2982.8333
$ws.Range("J122").Value = 2966.6667
$ws.Range("L122").Value = 8900.000100000001
$ws.Range("N122").Value = -13800.0001

$ws.Range("H136").Value = 2333.3333
$ws.Range("J136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496

$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H122").Value = 832.9583
$ws.Range("I122").Value = 851.7826
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 2555.3478
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = -105.3478
$ws.Range("N122").Value = -6100

$ws.Range("H126").Value = 6557.143
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530
